$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Page title (Heading1) and the later bold "title" run further down the
#    page share the exact same original text, so a single global
#    Find/Replace takes care of both occurrences.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Play Fairy Dust Extreme Free - Review of Fantasy-themed Slot",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Play Fairy Dust Extreme for Free", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. "What we like" bullet list: reorder the items and retitle one of them.
#    Before: Impressive graphics and animations / Unique bonus features /
#            Intuitive and user-friendly interface / Decent RTP
#    After:  Intuitive and user-friendly interface / Impressive graphics and
#            animations / Exciting bonus features / Decent RTP
#    We rewrite the whole 4-paragraph block via InsertXML so the exact
#    paragraph/run layout (including the leading empty <w:r/> runs) is
#    reproduced faithfully.
# ---------------------------------------------------------------------------
$bulletStart = $null
$bulletEnd = $null
foreach ($p in $d.Paragraphs) {
    $styleName = $p.Range.ParagraphStyle.NameLocal
    if ($styleName -eq "List Bullet" -and $p.Range.Text -like "Impressive graphics*") {
        $bulletStart = $p
    }
    if ($styleName -eq "List Bullet" -and $p.Range.Text -like "Decent RTP*" -and $bulletStart -ne $null -and $bulletEnd -eq $null) {
        $bulletEnd = $p
    }
}

$rng = $d.Range($bulletStart.Range.Start, $bulletEnd.Range.End)
$nsW = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$pPr = '<w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr>'
$xml = "<w:p $nsW>$pPr<w:r/><w:r><w:t>Intuitive and user-friendly interface</w:t></w:r></w:p>" +
       "<w:p $nsW>$pPr<w:r/><w:r><w:t>Impressive graphics and animations</w:t></w:r></w:p>" +
       "<w:p $nsW>$pPr<w:r/><w:r><w:t>Exciting bonus features</w:t></w:r></w:p>" +
       "<w:p $nsW>$pPr<w:r/><w:r><w:t>Decent RTP</w:t></w:r></w:p>"
$rng.InsertXML($xml) | Out-Null

# ---------------------------------------------------------------------------
# 3. "What we don't like" bullet: "Limited number of pay lines" becomes
#    "Limited betting limits". Using InsertXML (rather than Find/Replace)
#    keeps the paragraph's leading empty <w:r/> run intact. Note: a
#    single-paragraph-for-single-paragraph InsertXML replacement
#    automatically preserves the target paragraph's own leading empty run,
#    so it must NOT be duplicated here.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $styleName = $p.Range.ParagraphStyle.NameLocal
    if ($styleName -eq "List Bullet" -and $p.Range.Text -like "Limited number of pay lines*") {
        $rng2 = $d.Range($p.Range.Start, $p.Range.End)
        $xml2 = "<w:p $nsW>$pPr<w:r><w:t>Limited betting limits</w:t></w:r></w:p>"
        $rng2.InsertXML($xml2) | Out-Null
        break
    }
}

# ---------------------------------------------------------------------------
# 4. Meta description (italic run near the end of the document).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Play Fairy Dust Extreme for free and read our review of this fantasy-themed online slot game. Find out the pros and cons of Fairy Dust Extreme.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Experience the excitement of Fairy Dust Extreme, a fantasy-themed online slot game. Play for free now!", 2) | Out-Null
